$wb = $excel.ActiveWorkbook

$oldGuidHash = "eedb76f8-181b-4934-a578-f1cc15b07225"
$newGuidHash = "0aa75caa-44d8-456c-9b37-6204ba854d46"
$sourceHash  = "1f41b0d89e588dd7099e52de5bbf137618219458"

$oldFileName = "$oldGuidHash.md"
$newFileName = "$newGuidHash.md"
$oldPathName = "e2e\$oldGuidHash.md"
$newPathName = "e2e\$newGuidHash.md"

$newHoDate       = "2016-08-16 22:56:19"
$newZhHandoffFile = "$newGuidHash.$sourceHash.zh-cn.xlf"
$newZhHandoffDate = "2016-08-16 22:56:15"
$newDeHandoffFile = "$newGuidHash.$sourceHash.de-de.xlf"

$blobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cbf8563a4319471b1746f6f093c1cf62e0a5d49/e2e/"

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$blobBase$newFileName", "", "", $newPathName)

# ---- Sheet "zh-cn" ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("G2").Value = $newZhHandoffFile
$wsZh.Range("H2").Value = $newZhHandoffDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$blobBase$newFileName", "", "", $newFileName)

# ---- Sheet "de-de" ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("G2").Value = $newDeHandoffFile
$wsDe.Range("H2").Value = $newHoDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$blobBase$newFileName", "", "", $newFileName)
